# "Elimina dato 2023 de poblacion censada"
# Remove the 2023 "Censos" observation (row 7: Censos / 2023 / 3444.263)
# and keep the rest of the auxiliary projections table intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row containing the 2023 census datum; rows below shift up.
$ws.Rows.Item(7).Delete()

# The data now spans A1:C14 (13 data rows + header). Re-apply the sort that
# was already in effect so the sheet's remembered sort range shrinks from
# A2:C14 to A2:C13, matching the reduced data set.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A13"))
$ws.Sort.SortFields.Add($ws.Range("B2:B13"))
$ws.Sort.SetRange($ws.Range("A2:C13"))
$ws.Sort.Header = 2
$ws.Sort.Apply()

# Shrink the (hidden, leftover) _FilterDatabase defined name to match the
# new data extent.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Hoja1!`$A`$1:`$C`$13"
    }
}

# Restore the last active-cell selection recorded in the sheet view.
$ws.Range("D17").Select()
